$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 7) for June 6th, matching the layout of the
# existing rows (A: index, B: date, C-F: counts, G: percentage).
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A7").Value = 5

$ws.Range("B7").Value = 43988
$ws.Range("B7").NumberFormat = "YYYY-MM-DD"

$ws.Range("C7").Value = 113619
$ws.Range("D7").Value = 170434
$ws.Range("E7").Value = 48273
$ws.Range("F7").Value = 13511
$ws.Range("G7").Value = 33.72
